$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (data) and D (id_venda) hold numeric-looking text (dates / ids) that must
# stay as text, not be auto-converted to a date serial / number by Excel. Force the
# cells to Text format before writing so the stored type matches the source data.
$ws.Range("A2:A9").NumberFormat = "@"
$ws.Range("D2:D9").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2,1).Value = "2025-06-19"
$ws.Cells.Item(2,4).Value = "365782"
$ws.Cells.Item(2,7).Value = -449

# Row 3
$ws.Cells.Item(3,1).Value = "2025-06-19"
$ws.Cells.Item(3,4).Value = "366707"
$ws.Cells.Item(3,5).Value = 13079
$ws.Cells.Item(3,6).Value = "FONE BLUETOOTH BASIKE TWS FON6694"
$ws.Cells.Item(3,7).Value = -449
$ws.Cells.Item(3,8).Value = 1.08

# Row 4
$ws.Cells.Item(4,1).Value = "2025-06-24"
$ws.Cells.Item(4,4).Value = "370495"
$ws.Cells.Item(4,5).Value = 46217
$ws.Cells.Item(4,6).Value = "SMART WATCH HMASTON INK12"
$ws.Cells.Item(4,7).Value = -88
$ws.Cells.Item(4,8).Value = 1.03
$ws.Cells.Item(4,9).Value = 0.18

# Row 5
$ws.Cells.Item(5,1).Value = "2025-06-30"
$ws.Cells.Item(5,4).Value = "374463"
$ws.Cells.Item(5,7).Value = -449

# Row 6
$ws.Cells.Item(6,1).Value = "2025-06-30"
$ws.Cells.Item(6,4).Value = "374471"
$ws.Cells.Item(6,5).Value = 49904
$ws.Cells.Item(6,6).Value = "FONE SEM FIO OWS EAR-HOOK BRANCO BASIKE"
$ws.Cells.Item(6,7).Value = -16
$ws.Cells.Item(6,8).Value = 1.07
$ws.Cells.Item(6,9).Value = 0.26

# Row 7
$ws.Cells.Item(7,4).Value = "374491"
$ws.Cells.Item(7,5).Value = 14186
$ws.Cells.Item(7,6).Value = "BALANÇA DIGITAL 10KG"
$ws.Cells.Item(7,7).Value = -133
$ws.Cells.Item(7,8).Value = 1.04
$ws.Cells.Item(7,9).Value = 0.19

# Row 8 (new row, was previously part of 2025-06-30/374471 entry, now 2025-07-01)
$ws.Cells.Item(8,1).Value = "2025-07-01"
$ws.Cells.Item(8,4).Value = "375040"
$ws.Cells.Item(8,5).Value = 10114
$ws.Cells.Item(8,6).Value = "CARREGADOR USB-C A GOLD 20W CA31-4"
$ws.Cells.Item(8,7).Value = -84
$ws.Cells.Item(8,8).Value = 1.05
$ws.Cells.Item(8,9).Value = 0.22

# Row 9 (new row, was previously 374491, now 2025-07-01 / 375082)
$ws.Cells.Item(9,1).Value = "2025-07-01"
$ws.Cells.Item(9,4).Value = "375082"
$ws.Cells.Item(9,5).Value = 13546
$ws.Cells.Item(9,6).Value = "FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON"
$ws.Cells.Item(9,7).Value = -312
$ws.Cells.Item(9,8).Value = 1.09
$ws.Cells.Item(9,9).Value = 0.29
